$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.758.51"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.594.76"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "209.32"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").Value = "0.501"
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("E10").Value = "  -1.55%  "
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("D12").Value = "1.821.70"
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("D13").Value = "1.597.55"
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("E15").Value = "  -3.38%  "
$ws.Range("D16").Value = "27.758.35"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "63.45"
$ws.Range("E17").Value = "  -1.51%  "
$ws.Range("D18").Value = "219.58"
$ws.Range("E18").Value = "  -2.81%  "
$ws.Range("E19").Value = "  -2.50%  "
$ws.Range("E20").Value = "  -1.99%  "
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("E22").Value = "  -3.03%  "
$ws.Range("E23").Value = "  -2.26%  "
$ws.Range("E24").Value = "  -3.57%  "
$ws.Range("D25").Value = "154.05"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("D26").Value = "7.11"
$ws.Range("E26").Value = "  +3.37%  "
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("E29").Value = "  -3.04%  "
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("E32").Value = "  -3.85%  "
$ws.Range("D33").Value = "1.381.00"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("D34").Value = "2.97"
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("E35").Value = "  -3.50%  "
$ws.Range("D36").Value = "0.975"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("D39").Value = "0.536"
$ws.Range("E39").Value = "  -2.79%  "
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").Value = "0.977"
$ws.Range("E42").Value = "  -1.97%  "
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("D47").Value = "1.732.46"
$ws.Range("E47").Value = "  -1.46%  "
$ws.Range("D48").Value = "86.27"
$ws.Range("E48").Value = "  -3.60%  "
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("D50").Value = "0.0967"
$ws.Range("E50").Value = "  -2.33%  "
$ws.Range("E51").Value = "  -1.06%  "
